$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "ADD: COLLECTION" (right-hand) table rows 3-5 ---
# Row 3: JADE MAJID -> ALEJANDRO PIZON, S.I. 4241 -> 4320, amount 785 -> 3925
$ws.Range("H3").Value = "ALEJANDRO PIZON"
$ws.Range("I3").Value = 4320
$ws.Range("J3").Value = 3925

# Row 4: MYRA APILAN -> MARY JANE LIWASAG, S.I. 4341 -> 4319, amount 785 -> 2355
$ws.Range("H4").Value = "MARY JANE LIWASAG"
$ws.Range("I4").Value = 4319
$ws.Range("J4").Value = 2355

# Row 5: JAY KIN IDIAS / 4344 / 2355 entry removed entirely
$ws.Range("H5").Value = $null
$ws.Range("I5").Value = $null
$ws.Range("J5").Value = $null

# --- Update the workbook's saved Print Area ---
$ws.PageSetup.PrintArea = '$H$1:$M$12'

# --- Update view: zoom level and current selection ---
$excel.ActiveWindow.Zoom = 115
[void]$ws.Range("J5:L5").Select()
